$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = @'
store/favorite
'@

$ws.Range("C7").ClearContents() | Out-Null

$ws.Range("D7").Value = @'
scroll
'@

$ws.Range("E7").ClearContents() | Out-Null

$ws.Range("F7").Value = @'
channel, page_url, scroll_rate, os_name
'@

$ws.Range("G7").Value = @'
Rround, https://store.rround.com/favorite?fromMypage=true&tab=goods, 50, iOS
'@

$ws.Range("H7").Value = 4

$ws.Range("E11").Value = @'
바비리스 버터 바 스트레이트너 ST520K
'@

$ws.Range("G11").Value = @'
Rround, https://store.rround.com/main/home, 바비리스 버터 바 스트레이트너 ST520K, 식품
, 5, 1/7, 8038, 바비리스 버터 바 스트레이트너 ST520K, 최우수판매대리점, 59,000원, 26,000원, 55%, 무료배송, F, 5, iOS
'@

$ws.Range("F12").Value = @'
channel, page_url, area_name, tab_name, prd_order, area_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, prd_is_ad, el_order, os_name
'@

$ws.Range("G12").Value = @'
Rround, https://store.rround.com/main/home, 주목할 만한 상품이에요!, 식품
, 6, 3/7, 636, 비스카 블루투스 스마트 체중계 VK-S2(블랙), 빅픽처코퍼레이션, 26,900원, 13,900원, 48%, 1, 4, 무료배송, F, 6, iOS
'@

$ws.Range("H12").Value = 18

$ws.Range("E15").Value = @'
생활
'@

$ws.Range("F15").Value = @'
channel, page_url, click_text, srch_kwd, os_name
'@

$ws.Range("G15").Value = @'
Rround, https://store.rround.com/main/ranking, 생활, 생활, iOS
'@

$ws.Range("H15").Value = 5

$ws.Range("E18").Value = @'
(답이답이다) 베이킹소다 액체 세탁세제 용기 3L 4개
'@

$ws.Range("F18").Value = @'
channel, page_url, click_text, tab_name, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, prd_is_ad, os_name
'@

$ws.Range("G18").Value = @'
Rround, https://store.rround.com/main/ranking, (답이답이다) 베이킹소다 액체 세탁세제 용기 3L 4개, 식품
, 2, 28, (답이답이다) 베이킹소다 액체 세탁세제 용기 3L 4개, 케이디글로벌, 16,900원, 15,900원, 5%, 32, 4.5, 무료배송, F, iOS
'@

$ws.Range("H18").Value = 16

$ws.Range("F19").Value = @'
channel, page_url, tab_name, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, prd_is_ad, os_name
'@

$ws.Range("G19").Value = @'
Rround, https://store.rround.com/main/ranking, 식품
, 3, 36, (답이답이다) 제습제 520ml 12개, 케이디글로벌, 12,900원, 12,200원, 5%, 3, 5, 무료배송, F, iOS
'@

$ws.Range("H19").Value = 15

$ws.Range("F22").Value = @'
channel, page_url, click_text, area_name, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_disc_rate, prd_is_ad, os_name
'@

$ws.Range("G22").Value = @'
Rround, https://store.rround.com/main/deal, [닥터포헤어] 1+1 바이오3 탈모완화 샴푸 500ml, 무더운 여름, 케어의 시작, 2, 8070, [닥터포헤어] 1+1 바이오3 탈모완화 샴푸 500ml, 닥터포헤어, 28,900원, 9%, F, iOS
'@

$ws.Range("H22").Value = 12

$ws.Range("F23").Value = @'
channel, page_url, area_name, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_disc_rate, prd_is_ad, os_name
'@

$ws.Range("G23").Value = @'
Rround, https://store.rround.com/main/deal, 무더운 여름, 케어의 시작, 3, 8067, [닥터지] 레드 블레미쉬 클리어 수딩토너 기획세트, 닥터지, 18,500원, 53%, F, iOS
'@

$ws.Range("H23").Value = 11

$ws.Range("B31").Value = @'
store/product/detail/8038
'@

$ws.Range("G31").Value = @'
Rround, https://store.rround.com/product/detail/8038, 구매하기, 상품상세
, 8038, 바비리스 버터 바 스트레이트너 ST520K, 59,000원, 26,000원, 55%, 0, 0, #고데기___#고대기___#여행용고데기___#미용실고데기___#뿌리볼륨고데기___#가벼운고데기___#스트레이트너___#매직기___#바비리스고데기___#웨이브고데기___#온도조절고데기, iOS
'@

$ws.Range("B32").Value = @'
store/product/detail/8038
'@

$ws.Range("C32").Value = @'
구매 버튼
'@

$ws.Range("D32").Value = @'
click
'@

$ws.Range("E32").Value = @'
바로 구매하기
'@

$ws.Range("F32").Value = @'
channel, page_url, click_text, tab_name, prd_code, prd_name, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name
'@

$ws.Range("G32").Value = @'
Rround, https://store.rround.com/product/detail/8038, 바로 구매하기, 상품상세
, 8038, 바비리스 버터 바 스트레이트너 ST520K, 59,000원, 26,000원, 55%, 0, 0, #고데기___#고대기___#여행용고데기___#미용실고데기___#뿌리볼륨고데기___#가벼운고데기___#스트레이트너___#매직기___#바비리스고데기___#웨이브고데기___#온도조절고데기, iOS
'@
